$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 63
$ws.Range("C45").Value = 2
$ws.Range("D45").Value = 10
$ws.Range("E45").Value = 25
$ws.Range("F45").Value = 75
$ws.Range("G45").Value = 100
